# Apply the "change source structure (dataquery, admin-portal, bookstore-portal)"
# edit to the Back-end sheet of the project-plan workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Back-end")
$ws.Activate()

# Row 7: task name was renamed from "CRUD ..." wording to "Create ..." wording.
$ws.Range("C7").Value = "Create User & Role entity"

# Row 11: task name rewritten, and the previously-empty effort/date/status cells
# are now filled in (mirrors the Back-end/Front-end effort columns).
$ws.Range("C11").Value = "Create all Entities"
$ws.Range("D11").Value = "8h"
$ws.Range("E11").Value = "23/4"
$ws.Range("F11").Value = "24/4"
$ws.Range("G11").Value = "8h"
$ws.Range("H11").Value = "23/4"
$ws.Range("I11").Value = "24/4"
$ws.Range("J11").Value = "Done"

# Row 12: a brand-new task describing the source-structure change, with its
# effort/date/status/note cells filled in.
$ws.Range("K12").Value = "3 part: dataquery, admin, bookstore"
$ws.Range("D12").Value = "6h"
$ws.Range("C12").Value = "Change source structure"
$ws.Range("E12").Value = "24/4"
$ws.Range("F12").Value = "24/4"
$ws.Range("G12").Value = "5h"
$ws.Range("H12").Value = "24/4"
$ws.Range("I12").Value = "24/4"
$ws.Range("J12").Value = "Done"

# Column K (Note) is widened to fit the longer note text now present.
$ws.Columns.Item(11).ColumnWidth = 37

# Move the active selection to reflect where the author ended up editing.
$ws.Range("F23").Select()
